# Fixes to the readme.
#
# "How voice is usually processed?" -> "How is voice usually processed?"
# (TextBox 44 on slide 1, next to the "Where does MicGate affect this?" box)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 44")
$tr = $shp.TextFrame.TextRange

# Original runs: "How voice is " | "usually " | "processed?"
# Shrink the first run down to "How " (keeps its own rPr/smtClean intact).
$run1 = $tr.Characters(1, 13)
$run1.Text = "How "

# The old second run ("usually ") now starts right after "How ".
# Prefix it with "is voice " so it reads "is voice usually ".
$run2 = $tr.Characters(5, 8)
$run2.Text = "is voice " + $run2.Text
